$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "procedencia"

# Update header values to lowercase, unaccented versions
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "tipo"
$ws.Range("D1").Value = "ubicacion"
$ws.Range("E1").Value = "contacto"
$ws.Range("F1").Value = "telefono"
$ws.Range("G1").Value = "descripcion"
$ws.Range("H1").Value = "comentario"

# Strip the old bold/white-on-blue centered header styling back to plain default
$ws.Range("A1:H1").ClearFormats()
